# Fixed issue - two digits missing for compensation number.
# This adds a new "API_Environment" asset to the Assets sheet, and fills in
# the previously-empty "Completed_MailSubject" / "Completed_MailBody"
# constants on the Constants sheet. Also restores the on-screen selection
# (active cell) for each worksheet.

$wb = $excel.ActiveWorkbook
$wsSettings  = $wb.Worksheets.Item("Settings")
$wsConstants = $wb.Worksheets.Item("Constants")
$wsAssets    = $wb.Worksheets.Item("Assets")

# --- Data changes -----------------------------------------------------
# Order of assignment matters so new shared-string entries land in the
# same order as the target workbook:
#   1128_API_Environment, API_Environment, Robot Finished Processing,
#   Hi, Robot successfully finished Order Finalization & Claim Settlement. Regards, EC_JD Robot.

# Assets sheet: new row 26 - API_Environment asset
$wsAssets.Range("B26").Value = "1128_API_Environment"
$wsAssets.Range("A26").Value = "API_Environment"

# Constants sheet: fill in the mail subject/body for the completed notification
$wsConstants.Range("B24").Value = "Robot Finished Processing"
$wsConstants.Range("B25").Value = "Hi, Robot successfully finished Order Finalization & Claim Settlement. Regards, EC_JD Robot."

# --- View state (selected cell per sheet) ------------------------------
# Activate sheets in order, leaving Assets active last (matches workbook's
# original activeTab / tabSelected state).
$wsSettings.Activate()
$wsSettings.Range("A5").Select()

$wsConstants.Activate()
$wsConstants.Range("B25").Select()

$wsAssets.Activate()
$wsAssets.Range("A7").Select()
